$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1 (Ghana/Kirguistan etc. reorder + daily case refresh)
$ws.Range("A1").Value = 'Datos actualizados a 27 de Septiembre de 2020 a las 17:31'

$data = New-Object 'object[,]' 216,8
$data[0,0] = 'Estados Unidos'
$data[0,1] = 7292796
$data[0,2] = 5235
$data[0,3] = 4524730
$data[0,4] = 2558836
$data[0,5] = 0
$data[0,6] = 53
$data[0,7] = 209230
$data[1,0] = 'India'
$data[1,1] = 6041638
$data[1,2] = 51057
$data[1,3] = 4981099
$data[1,4] = 965568
$data[1,5] = 0
$data[1,6] = 437
$data[1,7] = 94971
$data[2,0] = 'Brasil'
$data[2,1] = 4718115
$data[2,2] = 0
$data[2,3] = 4050837
$data[2,4] = 525837
$data[2,5] = 0
$data[2,6] = 0
$data[2,7] = 141441
$data[3,0] = 'Rusia'
$data[3,1] = 1151438
$data[3,2] = 7867
$data[3,3] = 943218
$data[3,4] = 187896
$data[3,5] = 0
$data[3,6] = 99
$data[3,7] = 20324
$data[4,0] = 'Colombia'
$data[4,1] = 806038
$data[4,2] = 0
$data[4,3] = 700112
$data[4,4] = 80630
$data[4,5] = 0
$data[4,6] = 0
$data[4,7] = 25296
$data[5,0] = 'Peru'
$data[5,1] = 800142
$data[5,2] = 0
$data[5,3] = 657836
$data[5,4] = 110164
$data[5,5] = 0
$data[5,6] = 0
$data[5,7] = 32142
$data[6,0] = 'España'
$data[6,1] = 735198
$data[6,2] = 0
$data[6,3] = 0
$data[6,4] = 0
$data[6,5] = 0
$data[6,6] = 0
$data[6,7] = 31232
$data[7,0] = 'Mexico'
$data[7,1] = 726431
$data[7,2] = 5573
$data[7,3] = 521241
$data[7,4] = 128947
$data[7,5] = 0
$data[7,6] = 399
$data[7,7] = 76243
$data[8,0] = 'Argentina'
$data[8,1] = 702484
$data[8,2] = 0
$data[8,3] = 556489
$data[8,4] = 130452
$data[8,5] = 0
$data[8,6] = 0
$data[8,7] = 15543
$data[9,0] = 'Sudafrica'
$data[9,1] = 669498
$data[9,2] = 0
$data[9,3] = 601818
$data[9,4] = 51304
$data[9,5] = 0
$data[9,6] = 0
$data[9,7] = 16376
$data[10,0] = 'Francia'
$data[10,1] = 527446
$data[10,2] = 0
$data[10,3] = 94891
$data[10,4] = 400855
$data[10,5] = 0
$data[10,6] = 0
$data[10,7] = 31700
$data[11,0] = 'Chile'
$data[11,1] = 457901
$data[11,2] = 1922
$data[11,3] = 431704
$data[11,4] = 13556
$data[11,5] = 0
$data[11,6] = 50
$data[11,7] = 12641
$data[12,0] = 'Iran'
$data[12,1] = 446448
$data[12,2] = 3362
$data[12,3] = 374170
$data[12,4] = 46689
$data[12,5] = 0
$data[12,6] = 195
$data[12,7] = 25589
$data[13,0] = 'Reino Unido'
$data[13,1] = 429277
$data[13,2] = 0
$data[13,3] = 0
$data[13,4] = 0
$data[13,5] = 0
$data[13,6] = 0
$data[13,7] = 41971
$data[14,0] = 'Banglades'
$data[14,1] = 359148
$data[14,2] = 1275
$data[14,3] = 270491
$data[14,4] = 83496
$data[14,5] = 0
$data[14,6] = 32
$data[14,7] = 5161
$data[15,0] = 'Irak'
$data[15,1] = 349450
$data[15,2] = 3481
$data[15,3] = 280673
$data[15,4] = 59787
$data[15,5] = 0
$data[15,6] = 55
$data[15,7] = 8990
$data[16,0] = 'Arabia Saudita'
$data[16,1] = 333193
$data[16,2] = 403
$data[16,3] = 317005
$data[16,4] = 11505
$data[16,5] = 0
$data[16,6] = 28
$data[16,7] = 4683
$data[17,0] = 'Turquia'
$data[17,1] = 312966
$data[17,2] = 0
$data[17,3] = 274514
$data[17,4] = 30523
$data[17,5] = 0
$data[17,6] = 0
$data[17,7] = 7929
$data[18,0] = 'Pakistan'
$data[18,1] = 310275
$data[18,2] = 694
$data[18,3] = 295613
$data[18,4] = 8205
$data[18,5] = 0
$data[18,6] = 6
$data[18,7] = 6457
$data[19,0] = 'Italia'
$data[19,1] = 309870
$data[19,2] = 1766
$data[19,3] = 224417
$data[19,4] = 49618
$data[19,5] = 0
$data[19,6] = 17
$data[19,7] = 35835
$data[20,0] = 'Filipinas'
$data[20,1] = 304226
$data[20,2] = 2995
$data[20,3] = 252510
$data[20,4] = 46372
$data[20,5] = 0
$data[20,6] = 60
$data[20,7] = 5344
$data[21,0] = 'Alemania'
$data[21,1] = 285270
$data[21,2] = 245
$data[21,3] = 249500
$data[21,4] = 26238
$data[21,5] = 0
$data[21,6] = 0
$data[21,7] = 9532
$data[22,0] = 'Indonesia'
$data[22,1] = 275213
$data[22,2] = 3874
$data[22,3] = 203014
$data[22,4] = 61813
$data[22,5] = 0
$data[22,6] = 78
$data[22,7] = 10386
$data[23,0] = 'Israel'
$data[23,1] = 231026
$data[23,2] = 3926
$data[23,3] = 159931
$data[23,4] = 69629
$data[23,5] = 0
$data[23,6] = 25
$data[23,7] = 1466
$data[24,0] = 'Ucrania'
$data[24,1] = 198634
$data[24,2] = 3130
$data[24,3] = 87882
$data[24,4] = 106793
$data[24,5] = 0
$data[24,6] = 56
$data[24,7] = 3959
$data[25,0] = 'Canada'
$data[25,1] = 152162
$data[25,2] = 491
$data[25,3] = 130617
$data[25,4] = 12281
$data[25,5] = 0
$data[25,6] = 2
$data[25,7] = 9264
$data[26,0] = 'Ecuador'
$data[26,1] = 133981
$data[26,2] = 0
$data[26,3] = 102852
$data[26,4] = 19856
$data[26,5] = 0
$data[26,6] = 0
$data[26,7] = 11273
$data[27,0] = 'Bolivia'
$data[27,1] = 133592
$data[27,2] = 370
$data[27,3] = 93406
$data[27,4] = 32358
$data[27,5] = 0
$data[27,6] = 28
$data[27,7] = 7828
$data[28,0] = 'Catar'
$data[28,1] = 125084
$data[28,2] = 234
$data[28,3] = 121995
$data[28,4] = 2875
$data[28,5] = 0
$data[28,6] = 0
$data[28,7] = 214
$data[29,0] = 'Rumania'
$data[29,1] = 122673
$data[29,2] = 1438
$data[29,3] = 98607
$data[29,4] = 19348
$data[29,5] = 0
$data[29,6] = 31
$data[29,7] = 4718
$data[30,0] = 'Marruecos'
$data[30,1] = 115241
$data[30,2] = 0
$data[30,3] = 94150
$data[30,4] = 19050
$data[30,5] = 0
$data[30,6] = 0
$data[30,7] = 2041
$data[31,0] = 'Belgica'
$data[31,1] = 112803
$data[31,2] = 1827
$data[31,3] = 19246
$data[31,4] = 83583
$data[31,5] = 0
$data[31,6] = 5
$data[31,7] = 9974
$data[32,0] = 'Paises Bajos'
$data[32,1] = 111626
$data[32,2] = 2995
$data[32,3] = 0
$data[32,4] = 0
$data[32,5] = 0
$data[32,6] = 8
$data[32,7] = 6374
$data[33,0] = 'Republica Dominicana'
$data[33,1] = 111386
$data[33,2] = 429
$data[33,3] = 85965
$data[33,4] = 23326
$data[33,5] = 0
$data[33,6] = 2
$data[33,7] = 2095
$data[34,0] = 'Panama'
$data[34,1] = 110108
$data[34,2] = 0
$data[34,3] = 86796
$data[34,4] = 20989
$data[34,5] = 0
$data[34,6] = 0
$data[34,7] = 2323
$data[35,0] = 'Kazajistan'
$data[35,1] = 107723
$data[35,2] = 64
$data[35,3] = 102666
$data[35,4] = 3358
$data[35,5] = 0
$data[35,6] = 0
$data[35,7] = 1699
$data[36,0] = 'Kuwait'
$data[36,1] = 103544
$data[36,2] = 345
$data[36,3] = 94929
$data[36,4] = 8014
$data[36,5] = 0
$data[36,6] = 4
$data[36,7] = 601
$data[37,0] = 'Egipto'
$data[37,1] = 102736
$data[37,2] = 0
$data[37,3] = 94374
$data[37,4] = 2493
$data[37,5] = 0
$data[37,6] = 0
$data[37,7] = 5869
$data[38,0] = 'Oman'
$data[38,1] = 97450
$data[38,2] = 1543
$data[38,3] = 87801
$data[38,4] = 8740
$data[38,5] = 0
$data[38,6] = 24
$data[38,7] = 909
$data[39,0] = 'Emiratos Arabes Unidos'
$data[39,1] = 91469
$data[39,2] = 851
$data[39,3] = 80544
$data[39,4] = 10513
$data[39,5] = 0
$data[39,6] = 1
$data[39,7] = 412
$data[40,0] = 'Suecia'
$data[40,1] = 90923
$data[40,2] = 0
$data[40,3] = 0
$data[40,4] = 0
$data[40,5] = 0
$data[40,6] = 0
$data[40,7] = 5880
$data[41,0] = 'Guatemala'
$data[41,1] = 90092
$data[41,2] = 390
$data[41,3] = 78698
$data[41,4] = 8165
$data[41,5] = 0
$data[41,6] = 16
$data[41,7] = 3229
$data[42,0] = 'Polonia'
$data[42,1] = 87330
$data[42,2] = 1350
$data[42,3] = 67904
$data[42,4] = 16994
$data[42,5] = 0
$data[42,6] = 8
$data[42,7] = 2432
$data[43,0] = 'China'
$data[43,1] = 85351
$data[43,2] = 14
$data[43,3] = 80541
$data[43,4] = 176
$data[43,5] = 0
$data[43,6] = 0
$data[43,7] = 4634
$data[44,0] = 'Japon'
$data[44,1] = 81055
$data[44,2] = 0
$data[44,3] = 74151
$data[44,4] = 5364
$data[44,5] = 0
$data[44,6] = 0
$data[44,7] = 1540
$data[45,0] = 'Bielorrusia'
$data[45,1] = 77609
$data[45,2] = 320
$data[45,3] = 74120
$data[45,4] = 2671
$data[45,5] = 0
$data[45,6] = 5
$data[45,7] = 818
$data[46,0] = 'Honduras'
$data[46,1] = 74548
$data[46,2] = 708
$data[46,3] = 26088
$data[46,4] = 46172
$data[46,5] = 0
$data[46,6] = 17
$data[46,7] = 2288
$data[47,0] = 'Portugal'
$data[47,1] = 73604
$data[47,2] = 665
$data[47,3] = 47647
$data[47,4] = 24004
$data[47,5] = 0
$data[47,6] = 9
$data[47,7] = 1953
$data[48,0] = 'Nepal'
$data[48,1] = 73394
$data[48,2] = 1573
$data[48,3] = 53898
$data[48,4] = 19019
$data[48,5] = 0
$data[48,6] = 10
$data[48,7] = 477
$data[49,0] = 'Etiopia'
$data[49,1] = 72700
$data[49,2] = 0
$data[49,3] = 30029
$data[49,4] = 41506
$data[49,5] = 0
$data[49,6] = 0
$data[49,7] = 1165
$data[50,0] = 'Costa Rica'
$data[50,1] = 72049
$data[50,2] = 0
$data[50,3] = 27760
$data[50,4] = 43461
$data[50,5] = 0
$data[50,6] = 0
$data[50,7] = 828
$data[51,0] = 'Venezuela'
$data[51,1] = 71940
$data[51,2] = 0
$data[51,3] = 61528
$data[51,4] = 9812
$data[51,5] = 0
$data[51,6] = 0
$data[51,7] = 600
$data[52,0] = 'Barein'
$data[52,1] = 68775
$data[52,2] = 0
$data[52,3] = 62252
$data[52,4] = 6281
$data[52,5] = 0
$data[52,6] = 3
$data[52,7] = 242
$data[53,0] = 'Chequia'
$data[53,1] = 63294
$data[53,2] = 0
$data[53,3] = 30936
$data[53,4] = 31767
$data[53,5] = 0
$data[53,6] = 0
$data[53,7] = 591
$data[54,0] = 'Nigeria'
$data[54,1] = 58198
$data[54,2] = 0
$data[54,3] = 49722
$data[54,4] = 7370
$data[54,5] = 0
$data[54,6] = 0
$data[54,7] = 1106
$data[55,0] = 'Singapur'
$data[55,1] = 57700
$data[55,2] = 15
$data[55,3] = 57367
$data[55,4] = 306
$data[55,5] = 0
$data[55,6] = 0
$data[55,7] = 27
$data[56,0] = 'Uzbekistan'
$data[56,1] = 55320
$data[56,2] = 501
$data[56,3] = 51829
$data[56,4] = 3033
$data[56,5] = 0
$data[56,6] = 6
$data[56,7] = 458
$data[57,0] = 'Suiza'
$data[57,1] = 51864
$data[57,2] = 0
$data[57,3] = 42600
$data[57,4] = 7200
$data[57,5] = 0
$data[57,6] = 0
$data[57,7] = 2064
$data[58,0] = 'Argelia'
$data[58,1] = 50914
$data[58,2] = 0
$data[58,3] = 35756
$data[58,4] = 13447
$data[58,5] = 0
$data[58,6] = 0
$data[58,7] = 1711
$data[59,0] = 'Moldavia'
$data[59,1] = 50875
$data[59,2] = 341
$data[59,3] = 37842
$data[59,4] = 11746
$data[59,5] = 0
$data[59,6] = 8
$data[59,7] = 1287
$data[60,0] = 'Armenia'
$data[60,1] = 49400
$data[60,2] = 328
$data[60,3] = 43613
$data[60,4] = 4836
$data[60,5] = 0
$data[60,6] = 3
$data[60,7] = 951
$data[61,0] = 'Ghana'
$data[61,1] = 46353
$data[61,2] = 131
$data[61,3] = 45577
$data[61,4] = 477
$data[61,5] = 0
$data[61,6] = 0
$data[61,7] = 299
$data[62,0] = 'Kirguistan'
$data[62,1] = 46251
$data[62,2] = 161
$data[62,3] = 42453
$data[62,4] = 2735
$data[62,5] = 0
$data[62,6] = 0
$data[62,7] = 1063
$data[63,0] = 'Austria'
$data[63,1] = 42876
$data[63,2] = 662
$data[63,3] = 33589
$data[63,4] = 8500
$data[63,5] = 0
$data[63,6] = 0
$data[63,7] = 787
$data[64,0] = 'Azerbaiyan'
$data[64,1] = 40023
$data[64,2] = 128
$data[64,3] = 37655
$data[64,4] = 1782
$data[64,5] = 0
$data[64,6] = 1
$data[64,7] = 586
$data[65,0] = 'Afganistan'
$data[65,1] = 39227
$data[65,2] = 35
$data[65,3] = 32642
$data[65,4] = 5132
$data[65,5] = 0
$data[65,6] = 0
$data[65,7] = 1453
$data[66,0] = 'Estado de Palestina'
$data[66,1] = 38703
$data[66,2] = 450
$data[66,3] = 29068
$data[66,4] = 9344
$data[66,5] = 0
$data[66,6] = 6
$data[66,7] = 291
$data[67,0] = 'Kenia'
$data[67,1] = 38115
$data[67,2] = 244
$data[67,3] = 24621
$data[67,4] = 12803
$data[67,5] = 0
$data[67,6] = 2
$data[67,7] = 691
$data[68,0] = 'Paraguay'
$data[68,1] = 37922
$data[68,2] = 0
$data[68,3] = 21757
$data[68,4] = 15383
$data[68,5] = 0
$data[68,6] = 0
$data[68,7] = 782
$data[69,0] = 'Libano'
$data[69,1] = 35242
$data[69,2] = 0
$data[69,3] = 15434
$data[69,4] = 19468
$data[69,5] = 0
$data[69,6] = 0
$data[69,7] = 340
$data[70,0] = 'Irlanda'
$data[70,1] = 34560
$data[70,2] = 0
$data[70,3] = 23364
$data[70,4] = 9394
$data[70,5] = 0
$data[70,6] = 0
$data[70,7] = 1802
$data[71,0] = 'Serbia'
$data[71,1] = 33384
$data[71,2] = 72
$data[71,3] = 31536
$data[71,4] = 1101
$data[71,5] = 0
$data[71,6] = 1
$data[71,7] = 747
$data[72,0] = 'Libia'
$data[72,1] = 32364
$data[72,2] = 536
$data[72,3] = 18128
$data[72,4] = 13716
$data[72,5] = 0
$data[72,6] = 21
$data[72,7] = 520
$data[73,0] = 'El Salvador'
$data[73,1] = 28630
$data[73,2] = 215
$data[73,3] = 22879
$data[73,4] = 4925
$data[73,5] = 0
$data[73,6] = 0
$data[73,7] = 826
$data[74,0] = 'Australia'
$data[74,1] = 27040
$data[74,2] = 24
$data[74,3] = 24573
$data[74,4] = 1595
$data[74,5] = 0
$data[74,6] = 2
$data[74,7] = 872
$data[75,0] = 'Bosnia y Herzegovina'
$data[75,1] = 26920
$data[75,2] = 123
$data[75,3] = 19746
$data[75,4] = 6352
$data[75,5] = 0
$data[75,6] = 2
$data[75,7] = 822
$data[76,0] = 'Dinamarca'
$data[76,1] = 26637
$data[76,2] = 424
$data[76,3] = 19650
$data[76,4] = 6338
$data[76,5] = 0
$data[76,6] = 1
$data[76,7] = 649
$data[77,0] = 'Hungria'
$data[77,1] = 24014
$data[77,2] = 937
$data[77,3] = 5141
$data[77,4] = 18137
$data[77,5] = 0
$data[77,6] = 6
$data[77,7] = 736
$data[78,0] = 'Corea del Sur'
$data[78,1] = 23611
$data[78,2] = 95
$data[78,3] = 21248
$data[78,4] = 1962
$data[78,5] = 0
$data[78,6] = 2
$data[78,7] = 401
$data[79,0] = 'Camerun'
$data[79,1] = 20735
$data[79,2] = 0
$data[79,3] = 19440
$data[79,4] = 877
$data[79,5] = 0
$data[79,6] = 0
$data[79,7] = 418
$data[80,0] = 'Bulgaria'
$data[80,1] = 19997
$data[80,2] = 0
$data[80,3] = 14160
$data[80,4] = 5048
$data[80,5] = 0
$data[80,6] = 0
$data[80,7] = 789
$data[81,0] = 'Costa de Marfil'
$data[81,1] = 19600
$data[81,2] = 0
$data[81,3] = 19122
$data[81,4] = 358
$data[81,5] = 0
$data[81,6] = 0
$data[81,7] = 120
$data[82,0] = 'Republica de Macedonia'
$data[82,1] = 17629
$data[82,2] = 146
$data[82,3] = 14581
$data[82,4] = 2323
$data[82,5] = 0
$data[82,6] = 3
$data[82,7] = 725
$data[83,0] = 'Grecia'
$data[83,1] = 17228
$data[83,2] = 0
$data[83,3] = 9989
$data[83,4] = 6863
$data[83,5] = 0
$data[83,6] = 0
$data[83,7] = 376
$data[84,0] = 'Madagascar'
$data[84,1] = 16285
$data[84,2] = 28
$data[84,3] = 14922
$data[84,4] = 1134
$data[84,5] = 0
$data[84,6] = 0
$data[84,7] = 229
$data[85,0] = 'Croacia'
$data[85,1] = 16197
$data[85,2] = 190
$data[85,3] = 14609
$data[85,4] = 1316
$data[85,5] = 0
$data[85,6] = 3
$data[85,7] = 272
$data[86,0] = 'Tunez'
$data[86,1] = 15178
$data[86,2] = 786
$data[86,3] = 5032
$data[86,4] = 9939
$data[86,5] = 0
$data[86,6] = 16
$data[86,7] = 207
$data[87,0] = 'Senegal'
$data[87,1] = 14909
$data[87,2] = 40
$data[87,3] = 12113
$data[87,4] = 2488
$data[87,5] = 0
$data[87,6] = 2
$data[87,7] = 308
$data[88,0] = 'Zambia'
$data[88,1] = 14612
$data[88,2] = 0
$data[88,3] = 13727
$data[88,4] = 553
$data[88,5] = 0
$data[88,6] = 0
$data[88,7] = 332
$data[89,0] = 'Noruega'
$data[89,1] = 13660
$data[89,2] = 33
$data[89,3] = 11190
$data[89,4] = 2200
$data[89,5] = 0
$data[89,6] = 0
$data[89,7] = 270
$data[90,0] = 'Sudan'
$data[90,1] = 13606
$data[90,2] = 0
$data[90,3] = 6764
$data[90,4] = 6006
$data[90,5] = 0
$data[90,6] = 0
$data[90,7] = 836
$data[91,0] = 'Albania'
$data[91,1] = 13259
$data[91,2] = 106
$data[91,3] = 7397
$data[91,4] = 5485
$data[91,5] = 0
$data[91,6] = 2
$data[91,7] = 377
$data[92,0] = 'Namibia'
$data[92,1] = 11033
$data[92,2] = 115
$data[92,3] = 8776
$data[92,4] = 2137
$data[92,5] = 0
$data[92,6] = 0
$data[92,7] = 120
$data[93,0] = 'Malasia'
$data[93,1] = 10919
$data[93,2] = 150
$data[93,3] = 9835
$data[93,4] = 950
$data[93,5] = 0
$data[93,6] = 1
$data[93,7] = 134
$data[94,0] = 'Birmania'
$data[94,1] = 10734
$data[94,2] = 743
$data[94,3] = 2862
$data[94,4] = 7646
$data[94,5] = 0
$data[94,6] = 28
$data[94,7] = 226
$data[95,0] = 'Consejo Danes para los Refugiados'
$data[95,1] = 10593
$data[95,2] = 0
$data[95,3] = 10093
$data[95,4] = 229
$data[95,5] = 0
$data[95,6] = 0
$data[95,7] = 271
$data[96,0] = 'Guinea'
$data[96,1] = 10512
$data[96,2] = 0
$data[96,3] = 9836
$data[96,4] = 611
$data[96,5] = 0
$data[96,6] = 0
$data[96,7] = 65
$data[97,0] = 'Montenegro'
$data[97,1] = 10197
$data[97,2] = 0
$data[97,3] = 6368
$data[97,4] = 3671
$data[97,5] = 0
$data[97,6] = 0
$data[97,7] = 158
$data[98,0] = 'Maldivas'
$data[98,1] = 10045
$data[98,2] = 0
$data[98,3] = 8754
$data[98,4] = 1257
$data[98,5] = 0
$data[98,6] = 0
$data[98,7] = 34
$data[99,0] = 'Guayana Francesa'
$data[99,1] = 9863
$data[99,2] = 0
$data[99,3] = 9500
$data[99,4] = 298
$data[99,5] = 0
$data[99,6] = 0
$data[99,7] = 65
$data[100,0] = 'Finlandia'
$data[100,1] = 9682
$data[100,2] = 105
$data[100,3] = 7850
$data[100,4] = 1489
$data[100,5] = 0
$data[100,6] = 0
$data[100,7] = 343
$data[101,0] = 'Tayikistan'
$data[101,1] = 9646
$data[101,2] = 41
$data[101,3] = 8430
$data[101,4] = 1141
$data[101,5] = 0
$data[101,6] = 0
$data[101,7] = 75
$data[102,0] = 'Eslovaquia'
$data[102,1] = 9078
$data[102,2] = 478
$data[102,3] = 4178
$data[102,4] = 4856
$data[102,5] = 0
$data[102,6] = 0
$data[102,7] = 44
$data[103,0] = 'Gabon'
$data[103,1] = 8728
$data[103,2] = 0
$data[103,3] = 7934
$data[103,4] = 740
$data[103,5] = 0
$data[103,6] = 0
$data[103,7] = 54
$data[104,0] = 'Haiti'
$data[104,1] = 8723
$data[104,2] = 0
$data[104,3] = 6551
$data[104,4] = 1945
$data[104,5] = 0
$data[104,6] = 0
$data[104,7] = 227
$data[105,0] = 'Jordania'
$data[105,1] = 8492
$data[105,2] = 431
$data[105,3] = 4222
$data[105,4] = 4225
$data[105,5] = 0
$data[105,6] = 2
$data[105,7] = 45
$data[106,0] = 'Luxemburgo'
$data[106,1] = 8311
$data[106,2] = 0
$data[106,3] = 6976
$data[106,4] = 1211
$data[106,5] = 0
$data[106,6] = 0
$data[106,7] = 124
$data[107,0] = 'Zimbabue'
$data[107,1] = 7803
$data[107,2] = 0
$data[107,3] = 6067
$data[107,4] = 1509
$data[107,5] = 0
$data[107,6] = 0
$data[107,7] = 227
$data[108,0] = 'Mozambique'
$data[108,1] = 7757
$data[108,2] = 0
$data[108,3] = 4769
$data[108,4] = 2934
$data[108,5] = 0
$data[108,6] = 0
$data[108,7] = 54
$data[109,0] = 'Uganda'
$data[109,1] = 7530
$data[109,2] = 166
$data[109,3] = 3647
$data[109,4] = 3810
$data[109,5] = 0
$data[109,6] = 2
$data[109,7] = 73
$data[110,0] = 'Mauritania'
$data[110,1] = 7462
$data[110,2] = 0
$data[110,3] = 7070
$data[110,4] = 231
$data[110,5] = 0
$data[110,6] = 0
$data[110,7] = 161
$data[111,0] = 'Jamaica'
$data[111,1] = 6017
$data[111,2] = 163
$data[111,3] = 1706
$data[111,4] = 4222
$data[111,5] = 0
$data[111,6] = 1
$data[111,7] = 89
$data[112,0] = 'Malaui'
$data[112,1] = 5766
$data[112,2] = 0
$data[112,3] = 4185
$data[112,4] = 1402
$data[112,5] = 0
$data[112,6] = 0
$data[112,7] = 179
$data[113,0] = 'Cabo Verde'
$data[113,1] = 5701
$data[113,2] = 0
$data[113,3] = 5018
$data[113,4] = 627
$data[113,5] = 0
$data[113,6] = 0
$data[113,7] = 56
$data[114,0] = 'Cuba'
$data[114,1] = 5457
$data[114,2] = 45
$data[114,3] = 4751
$data[114,4] = 584
$data[114,5] = 0
$data[114,6] = 2
$data[114,7] = 122
$data[115,0] = 'Suazilandia'
$data[115,1] = 5419
$data[115,2] = 0
$data[115,3] = 4802
$data[115,4] = 509
$data[115,5] = 0
$data[115,6] = 0
$data[115,7] = 108
$data[116,0] = 'Republica de Yibuti'
$data[116,1] = 5409
$data[116,2] = 0
$data[116,3] = 5340
$data[116,4] = 8
$data[116,5] = 0
$data[116,6] = 0
$data[116,7] = 61
$data[117,0] = 'Eslovenia'
$data[117,1] = 5350
$data[117,2] = 159
$data[117,3] = 3555
$data[117,4] = 1648
$data[117,5] = 0
$data[117,6] = 1
$data[117,7] = 147
$data[118,0] = 'Georgia'
$data[118,1] = 5254
$data[118,2] = 294
$data[118,3] = 1906
$data[118,4] = 3320
$data[118,5] = 0
$data[118,6] = 0
$data[118,7] = 28
$data[119,0] = 'Nicaragua'
$data[119,1] = 5073
$data[119,2] = 0
$data[119,3] = 2913
$data[119,4] = 2011
$data[119,5] = 0
$data[119,6] = 0
$data[119,7] = 149
$data[120,0] = 'Hong Kong'
$data[120,1] = 5066
$data[120,2] = 6
$data[120,3] = 4786
$data[120,4] = 175
$data[120,5] = 0
$data[120,6] = 0
$data[120,7] = 105
$data[121,0] = 'Guinea Ecuatorial'
$data[121,1] = 5028
$data[121,2] = 0
$data[121,3] = 4740
$data[121,4] = 205
$data[121,5] = 0
$data[121,6] = 0
$data[121,7] = 83
$data[122,0] = 'Congo'
$data[122,1] = 5005
$data[122,2] = 0
$data[122,3] = 3887
$data[122,4] = 1029
$data[122,5] = 0
$data[122,6] = 0
$data[122,7] = 89
$data[123,0] = 'Surinam'
$data[123,1] = 4831
$data[123,2] = 0
$data[123,3] = 4620
$data[123,4] = 109
$data[123,5] = 0
$data[123,6] = 0
$data[123,7] = 102
$data[124,0] = 'Ruanda'
$data[124,1] = 4811
$data[124,2] = 0
$data[124,3] = 3091
$data[124,4] = 1691
$data[124,5] = 0
$data[124,6] = 0
$data[124,7] = 29
$data[125,0] = 'Republica de Africa Central'
$data[125,1] = 4806
$data[125,2] = 0
$data[125,3] = 1840
$data[125,4] = 2904
$data[125,5] = 0
$data[125,6] = 0
$data[125,7] = 62
$data[126,0] = 'Angola'
$data[126,1] = 4672
$data[126,2] = 0
$data[126,3] = 1639
$data[126,4] = 2862
$data[126,5] = 0
$data[126,6] = 0
$data[126,7] = 171
$data[127,0] = 'Guadalupe'
$data[127,1] = 4487
$data[127,2] = 0
$data[127,3] = 2199
$data[127,4] = 2246
$data[127,5] = 0
$data[127,6] = 0
$data[127,7] = 42
$data[128,0] = 'Lituania'
$data[128,1] = 4385
$data[128,2] = 90
$data[128,3] = 2327
$data[128,4] = 1967
$data[128,5] = 0
$data[128,6] = 2
$data[128,7] = 91
$data[129,0] = 'Trinidad yTobago'
$data[129,1] = 4321
$data[129,2] = 9
$data[129,3] = 2234
$data[129,4] = 2017
$data[129,5] = 0
$data[129,6] = 0
$data[129,7] = 70
$data[130,0] = 'Siria'
$data[130,1] = 4038
$data[130,2] = 0
$data[130,3] = 1048
$data[130,4] = 2802
$data[130,5] = 0
$data[130,6] = 0
$data[130,7] = 188
$data[131,0] = 'Aruba'
$data[131,1] = 3832
$data[131,2] = 0
$data[131,3] = 2829
$data[131,4] = 978
$data[131,5] = 0
$data[131,6] = 0
$data[131,7] = 25
$data[132,0] = 'Bahamas'
$data[132,1] = 3790
$data[132,2] = 0
$data[132,3] = 1999
$data[132,4] = 1702
$data[132,5] = 0
$data[132,6] = 0
$data[132,7] = 89
$data[133,0] = 'Reunion'
$data[133,1] = 3685
$data[133,2] = 0
$data[133,3] = 2819
$data[133,4] = 855
$data[133,5] = 0
$data[133,6] = 0
$data[133,7] = 11
$data[134,0] = 'Somalia'
$data[134,1] = 3588
$data[134,2] = 0
$data[134,3] = 2943
$data[134,4] = 546
$data[134,5] = 0
$data[134,6] = 0
$data[134,7] = 99
$data[135,0] = 'Gambia'
$data[135,1] = 3564
$data[135,2] = 9
$data[135,3] = 2061
$data[135,4] = 1393
$data[135,5] = 0
$data[135,6] = 0
$data[135,7] = 110
$data[136,0] = 'Mayotte'
$data[136,1] = 3541
$data[136,2] = 0
$data[136,3] = 2964
$data[136,4] = 537
$data[136,5] = 0
$data[136,6] = 0
$data[136,7] = 40
$data[137,0] = 'Tailandia'
$data[137,1] = 3523
$data[137,2] = 1
$data[137,3] = 3367
$data[137,4] = 97
$data[137,5] = 0
$data[137,6] = 0
$data[137,7] = 59
$data[138,0] = 'Sri Lanka'
$data[138,1] = 3360
$data[138,2] = 11
$data[138,3] = 3208
$data[138,4] = 139
$data[138,5] = 0
$data[138,6] = 0
$data[138,7] = 13
$data[139,0] = 'Estonia'
$data[139,1] = 3200
$data[139,2] = 35
$data[139,3] = 2506
$data[139,4] = 630
$data[139,5] = 0
$data[139,6] = 0
$data[139,7] = 64
$data[140,0] = 'Mali'
$data[140,1] = 3080
$data[140,2] = 0
$data[140,3] = 2410
$data[140,4] = 540
$data[140,5] = 0
$data[140,6] = 0
$data[140,7] = 130
$data[141,0] = 'Malta'
$data[141,1] = 2979
$data[141,2] = 21
$data[141,3] = 2358
$data[141,4] = 590
$data[141,5] = 0
$data[141,6] = 0
$data[141,7] = 31
$data[142,0] = 'Botsuana'
$data[142,1] = 2921
$data[142,2] = 0
$data[142,3] = 701
$data[142,4] = 2204
$data[142,5] = 0
$data[142,6] = 0
$data[142,7] = 16
$data[143,0] = 'Guyana'
$data[143,1] = 2725
$data[143,2] = 0
$data[143,3] = 1535
$data[143,4] = 1116
$data[143,5] = 0
$data[143,6] = 0
$data[143,7] = 74
$data[144,0] = 'Sudan del Sur'
$data[144,1] = 2676
$data[144,2] = 0
$data[144,3] = 1290
$data[144,4] = 1337
$data[144,5] = 0
$data[144,6] = 0
$data[144,7] = 49
$data[145,0] = 'Islandia'
$data[145,1] = 2623
$data[145,2] = 22
$data[145,3] = 2158
$data[145,4] = 455
$data[145,5] = 0
$data[145,6] = 0
$data[145,7] = 10
$data[146,0] = 'Benin'
$data[146,1] = 2325
$data[146,2] = 0
$data[146,3] = 1960
$data[146,4] = 325
$data[146,5] = 0
$data[146,6] = 0
$data[146,7] = 40
$data[147,0] = 'Guinea-Bisau'
$data[147,1] = 2324
$data[147,2] = 0
$data[147,3] = 1549
$data[147,4] = 736
$data[147,5] = 0
$data[147,6] = 0
$data[147,7] = 39
$data[148,0] = 'Sierra Leona'
$data[148,1] = 2208
$data[148,2] = 0
$data[148,3] = 1679
$data[148,4] = 457
$data[148,5] = 0
$data[148,6] = 0
$data[148,7] = 72
$data[149,0] = 'Yemen'
$data[149,1] = 2030
$data[149,2] = 0
$data[149,3] = 1260
$data[149,4] = 183
$data[149,5] = 0
$data[149,6] = 0
$data[149,7] = 587
$data[150,0] = 'Uruguay'
$data[150,1] = 1998
$data[150,2] = 0
$data[150,3] = 1716
$data[150,4] = 235
$data[150,5] = 0
$data[150,6] = 0
$data[150,7] = 47
$data[151,0] = 'Burkina Faso'
$data[151,1] = 1973
$data[151,2] = 0
$data[151,3] = 1264
$data[151,4] = 653
$data[151,5] = 0
$data[151,6] = 0
$data[151,7] = 56
$data[152,0] = 'Principado de Andorra'
$data[152,1] = 1836
$data[152,2] = 0
$data[152,3] = 1263
$data[152,4] = 520
$data[152,5] = 0
$data[152,6] = 0
$data[152,7] = 53
$data[153,0] = 'Nueva Zelanda'
$data[153,1] = 1833
$data[153,2] = 2
$data[153,3] = 1749
$data[153,4] = 59
$data[153,5] = 0
$data[153,6] = 0
$data[153,7] = 25
$data[154,0] = 'Belice'
$data[154,1] = 1825
$data[154,2] = 17
$data[154,3] = 1165
$data[154,4] = 636
$data[154,5] = 0
$data[154,6] = 1
$data[154,7] = 24
$data[155,0] = 'Togo'
$data[155,1] = 1736
$data[155,2] = 0
$data[155,3] = 1319
$data[155,4] = 371
$data[155,5] = 0
$data[155,6] = 0
$data[155,7] = 46
$data[156,0] = 'Republica de Chipre'
$data[156,1] = 1684
$data[156,2] = 0
$data[156,3] = 1369
$data[156,4] = 293
$data[156,5] = 0
$data[156,6] = 0
$data[156,7] = 22
$data[157,0] = 'Letonia'
$data[157,1] = 1676
$data[157,2] = 22
$data[157,3] = 1304
$data[157,4] = 336
$data[157,5] = 0
$data[157,6] = 0
$data[157,7] = 36
$data[158,0] = 'Polinesia Francesa'
$data[158,1] = 1579
$data[158,2] = 0
$data[158,3] = 1335
$data[158,4] = 238
$data[158,5] = 0
$data[158,6] = 0
$data[158,7] = 6
$data[159,0] = 'Lesoto'
$data[159,1] = 1558
$data[159,2] = 0
$data[159,3] = 797
$data[159,4] = 726
$data[159,5] = 0
$data[159,6] = 0
$data[159,7] = 35
$data[160,0] = 'Liberia'
$data[160,1] = 1338
$data[160,2] = 0
$data[160,3] = 1221
$data[160,4] = 35
$data[160,5] = 0
$data[160,6] = 0
$data[160,7] = 82
$data[161,0] = 'Martinica'
$data[161,1] = 1290
$data[161,2] = 0
$data[161,3] = 98
$data[161,4] = 1172
$data[161,5] = 0
$data[161,6] = 0
$data[161,7] = 20
$data[162,0] = 'Niger'
$data[162,1] = 1194
$data[162,2] = 0
$data[162,3] = 1107
$data[162,4] = 18
$data[162,5] = 0
$data[162,6] = 0
$data[162,7] = 69
$data[163,0] = 'Republica del Chad'
$data[163,1] = 1177
$data[163,2] = 0
$data[163,3] = 1005
$data[163,4] = 89
$data[163,5] = 0
$data[163,6] = 0
$data[163,7] = 83
$data[164,0] = 'Vietnam'
$data[164,1] = 1074
$data[164,2] = 5
$data[164,3] = 999
$data[164,4] = 40
$data[164,5] = 0
$data[164,6] = 0
$data[164,7] = 35
$data[165,0] = 'Santo Tome y Principe'
$data[165,1] = 911
$data[165,2] = 0
$data[165,3] = 883
$data[165,4] = 13
$data[165,5] = 0
$data[165,6] = 0
$data[165,7] = 15
$data[166,0] = 'San Marino'
$data[166,1] = 727
$data[166,2] = 0
$data[166,3] = 676
$data[166,4] = 9
$data[166,5] = 0
$data[166,6] = 0
$data[166,7] = 42
$data[167,0] = 'Crucero'
$data[167,1] = 712
$data[167,2] = 0
$data[167,3] = 651
$data[167,4] = 48
$data[167,5] = 0
$data[167,6] = 0
$data[167,7] = 13
$data[168,0] = 'Islas Turcas y Caicos'
$data[168,1] = 681
$data[168,2] = 1
$data[168,3] = 620
$data[168,4] = 56
$data[168,5] = 0
$data[168,6] = 0
$data[168,7] = 5
$data[169,0] = 'San Martin (Parte Holandesa)'
$data[169,1] = 633
$data[169,2] = 6
$data[169,3] = 532
$data[169,4] = 79
$data[169,5] = 0
$data[169,6] = 0
$data[169,7] = 22
$data[170,0] = 'Papua Nueva Guinea'
$data[170,1] = 532
$data[170,2] = 0
$data[170,3] = 516
$data[170,4] = 9
$data[170,5] = 0
$data[170,6] = 0
$data[170,7] = 7
$data[171,0] = 'Taiwan'
$data[171,1] = 510
$data[171,2] = 0
$data[171,3] = 480
$data[171,4] = 23
$data[171,5] = 0
$data[171,6] = 0
$data[171,7] = 7
$data[172,0] = 'Tanzania'
$data[172,1] = 509
$data[172,2] = 0
$data[172,3] = 183
$data[172,4] = 305
$data[172,5] = 0
$data[172,6] = 0
$data[172,7] = 21
$data[173,0] = 'Burundi'
$data[173,1] = 485
$data[173,2] = 0
$data[173,3] = 472
$data[173,4] = 12
$data[173,5] = 0
$data[173,6] = 0
$data[173,7] = 1
$data[174,0] = 'Comoras'
$data[174,1] = 478
$data[174,2] = 0
$data[174,3] = 458
$data[174,4] = 13
$data[174,5] = 0
$data[174,6] = 0
$data[174,7] = 7
$data[175,0] = 'Islas Feroe'
$data[175,1] = 460
$data[175,2] = 0
$data[175,3] = 423
$data[175,4] = 37
$data[175,5] = 0
$data[175,6] = 0
$data[175,7] = 0
$data[176,0] = 'San Martin (Parte Francesa)'
$data[176,1] = 383
$data[176,2] = 16
$data[176,3] = 273
$data[176,4] = 102
$data[176,5] = 0
$data[176,6] = 0
$data[176,7] = 8
$data[177,0] = 'Gibraltar'
$data[177,1] = 379
$data[177,2] = 7
$data[177,3] = 341
$data[177,4] = 38
$data[177,5] = 0
$data[177,6] = 0
$data[177,7] = 0
$data[178,0] = 'Eritrea'
$data[178,1] = 375
$data[178,2] = 0
$data[178,3] = 341
$data[178,4] = 34
$data[178,5] = 0
$data[178,6] = 0
$data[178,7] = 0
$data[179,0] = 'Mauricio'
$data[179,1] = 367
$data[179,2] = 0
$data[179,3] = 343
$data[179,4] = 14
$data[179,5] = 0
$data[179,6] = 0
$data[179,7] = 10
$data[180,0] = 'Isla de Man'
$data[180,1] = 340
$data[180,2] = 0
$data[180,3] = 314
$data[180,4] = 2
$data[180,5] = 0
$data[180,6] = 0
$data[180,7] = 24
$data[181,0] = 'Curazao'
$data[181,1] = 337
$data[181,2] = 0
$data[181,3] = 134
$data[181,4] = 202
$data[181,5] = 0
$data[181,6] = 0
$data[181,7] = 1
$data[182,0] = 'Mongolia'
$data[182,1] = 313
$data[182,2] = 0
$data[182,3] = 303
$data[182,4] = 10
$data[182,5] = 0
$data[182,6] = 0
$data[182,7] = 0
$data[183,0] = 'Camboya'
$data[183,1] = 276
$data[183,2] = 1
$data[183,3] = 274
$data[183,4] = 2
$data[183,5] = 0
$data[183,6] = 0
$data[183,7] = 0
$data[184,0] = 'Butan'
$data[184,1] = 271
$data[184,2] = 8
$data[184,3] = 205
$data[184,4] = 66
$data[184,5] = 0
$data[184,6] = 0
$data[184,7] = 0
$data[185,0] = 'Monaco'
$data[185,1] = 210
$data[185,2] = 0
$data[185,3] = 169
$data[185,4] = 39
$data[185,5] = 0
$data[185,6] = 0
$data[185,7] = 2
$data[186,0] = 'Islas Caimanes'
$data[186,1] = 210
$data[186,2] = 0
$data[186,3] = 207
$data[186,4] = 2
$data[186,5] = 0
$data[186,6] = 0
$data[186,7] = 1
$data[187,0] = 'Barbados'
$data[187,1] = 190
$data[187,2] = 0
$data[187,3] = 178
$data[187,4] = 5
$data[187,5] = 0
$data[187,6] = 0
$data[187,7] = 7
$data[188,0] = 'Bermudas'
$data[188,1] = 181
$data[188,2] = 0
$data[188,3] = 167
$data[188,4] = 5
$data[188,5] = 0
$data[188,6] = 0
$data[188,7] = 9
$data[189,0] = 'Brunei'
$data[189,1] = 146
$data[189,2] = 0
$data[189,3] = 142
$data[189,4] = 1
$data[189,5] = 0
$data[189,6] = 0
$data[189,7] = 3
$data[190,0] = 'Seychelles'
$data[190,1] = 143
$data[190,2] = 0
$data[190,3] = 140
$data[190,4] = 3
$data[190,5] = 0
$data[190,6] = 0
$data[190,7] = 0
$data[191,0] = 'Liechtenstein'
$data[191,1] = 117
$data[191,2] = 0
$data[191,3] = 110
$data[191,4] = 6
$data[191,5] = 0
$data[191,6] = 0
$data[191,7] = 1
$data[192,0] = 'Antigua y Barbuda'
$data[192,1] = 98
$data[192,2] = 0
$data[192,3] = 92
$data[192,4] = 3
$data[192,5] = 0
$data[192,6] = 0
$data[192,7] = 3
$data[193,0] = 'Bonaire, San Eustaquio y Saba'
$data[193,1] = 85
$data[193,2] = 0
$data[193,3] = 21
$data[193,4] = 63
$data[193,5] = 0
$data[193,6] = 0
$data[193,7] = 1
$data[194,0] = 'Islas Virgenes Britanicas'
$data[194,1] = 71
$data[194,2] = 0
$data[194,3] = 62
$data[194,4] = 8
$data[194,5] = 0
$data[194,6] = 0
$data[194,7] = 1
$data[195,0] = 'San Vicente y las Granadinas'
$data[195,1] = 64
$data[195,2] = 0
$data[195,3] = 64
$data[195,4] = 0
$data[195,5] = 0
$data[195,6] = 0
$data[195,7] = 0
$data[196,0] = 'San Bartolome'
$data[196,1] = 48
$data[196,2] = 3
$data[196,3] = 25
$data[196,4] = 23
$data[196,5] = 0
$data[196,6] = 0
$data[196,7] = 0
$data[197,0] = 'Macao'
$data[197,1] = 46
$data[197,2] = 0
$data[197,3] = 46
$data[197,4] = 0
$data[197,5] = 0
$data[197,6] = 0
$data[197,7] = 0
$data[198,0] = 'Puerto Rico'
$data[198,1] = 39
$data[198,2] = 0
$data[198,3] = 1
$data[198,4] = 36
$data[198,5] = 0
$data[198,6] = 0
$data[198,7] = 2
$data[199,0] = 'Guam'
$data[199,1] = 32
$data[199,2] = 0
$data[199,3] = 0
$data[199,4] = 31
$data[199,5] = 0
$data[199,6] = 0
$data[199,7] = 1
$data[200,0] = 'Fiyi'
$data[200,1] = 32
$data[200,2] = 0
$data[200,3] = 28
$data[200,4] = 2
$data[200,5] = 0
$data[200,6] = 0
$data[200,7] = 2
$data[201,0] = 'Dominica'
$data[201,1] = 30
$data[201,2] = 6
$data[201,3] = 24
$data[201,4] = 6
$data[201,5] = 0
$data[201,6] = 0
$data[201,7] = 0
$data[202,0] = 'Nueva Caledonia'
$data[202,1] = 27
$data[202,2] = 0
$data[202,3] = 26
$data[202,4] = 1
$data[202,5] = 0
$data[202,6] = 0
$data[202,7] = 0
$data[203,0] = 'Santa Lucia'
$data[203,1] = 27
$data[203,2] = 0
$data[203,3] = 27
$data[203,4] = 0
$data[203,5] = 0
$data[203,6] = 0
$data[203,7] = 0
$data[204,0] = 'Timor Oriental'
$data[204,1] = 27
$data[204,2] = 0
$data[204,3] = 27
$data[204,4] = 0
$data[204,5] = 0
$data[204,6] = 0
$data[204,7] = 0
$data[205,0] = 'Granada'
$data[205,1] = 24
$data[205,2] = 0
$data[205,3] = 24
$data[205,4] = 0
$data[205,5] = 0
$data[205,6] = 0
$data[205,7] = 0
$data[206,0] = 'Laos'
$data[206,1] = 23
$data[206,2] = 0
$data[206,3] = 22
$data[206,4] = 1
$data[206,5] = 0
$data[206,6] = 0
$data[206,7] = 0
$data[207,0] = 'San Cristobal y Nieves'
$data[207,1] = 19
$data[207,2] = 0
$data[207,3] = 17
$data[207,4] = 2
$data[207,5] = 0
$data[207,6] = 0
$data[207,7] = 0
$data[208,0] = 'Islas Virgenes de los Estados Unidos'
$data[208,1] = 17
$data[208,2] = 0
$data[208,3] = 0
$data[208,4] = 17
$data[208,5] = 0
$data[208,6] = 0
$data[208,7] = 0
$data[209,0] = 'San Pedro y Miquelon'
$data[209,1] = 16
$data[209,2] = 0
$data[209,3] = 6
$data[209,4] = 10
$data[209,5] = 0
$data[209,6] = 0
$data[209,7] = 0
$data[210,0] = 'Groenlandia'
$data[210,1] = 14
$data[210,2] = 0
$data[210,3] = 14
$data[210,4] = 0
$data[210,5] = 0
$data[210,6] = 0
$data[210,7] = 0
$data[211,0] = 'Montserrat'
$data[211,1] = 13
$data[211,2] = 0
$data[211,3] = 12
$data[211,4] = 0
$data[211,5] = 0
$data[211,6] = 0
$data[211,7] = 1
$data[212,0] = 'Islas Malvinas'
$data[212,1] = 13
$data[212,2] = 0
$data[212,3] = 13
$data[212,4] = 0
$data[212,5] = 0
$data[212,6] = 0
$data[212,7] = 0
$data[213,0] = 'Santa Sede'
$data[213,1] = 12
$data[213,2] = 0
$data[213,3] = 12
$data[213,4] = 0
$data[213,5] = 0
$data[213,6] = 0
$data[213,7] = 0
$data[214,0] = 'Sahara Occidental'
$data[214,1] = 10
$data[214,2] = 0
$data[214,3] = 8
$data[214,4] = 1
$data[214,5] = 0
$data[214,6] = 0
$data[214,7] = 1
$data[215,0] = 'Anguila'
$data[215,1] = 3
$data[215,2] = 0
$data[215,3] = 3
$data[215,4] = 0
$data[215,5] = 0
$data[215,6] = 0
$data[215,7] = 0

$ws.Range("A4:H219").Value = $data

Write-Output "Updated countries table + timestamp"
